$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 438; existing rows 438-532 shift down to 439-533.
$ws.Rows(438).Insert()

# Populate the newly inserted row with the new price-report record.
$ws.Range("A438").Value = 10
$ws.Range("B438").Value = "Vega Modelo de Temuco"
$ws.Range("C438").Value = "La Araucanía"
$ws.Range("D438").Value = 44798
$ws.Range("E438").Value = 9
$ws.Range("F438").Value = "Fruta"
$ws.Range("G438").Value = 100108
$ws.Range("H438").Value = "Tropicales y subtropicales"
$ws.Range("I438").Value = 100108005
$ws.Range("J438").Value = "Piña"
$ws.Range("K438").Value = "Caramelo"
$ws.Range("L438").Value = "Segunda"
$ws.Range("M438").Value = 120
$ws.Range("N438").Value = 20000
$ws.Range("O438").Value = 22000
$ws.Range("P438").Value = 20917
$ws.Range("Q438").Value = "$/caja 14 unidades"
$ws.Range("R438").Value = "Ecuador"
$ws.Range("S438").Value = 1494
$ws.Range("T438").Value = 14
